$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text casing (values stay semantically the same, just lowercased)
# Order of assignment controls the shared-string table insertion order.
$ws.Range("A2").Value = "irrigation_volume"
$ws.Range("A3").Value = "revenue"
$ws.Range("D1").Value = "maize_cassava_beans"
$ws.Range("C1").Value = "vegetables"
$ws.Range("B1").Value = "fruits"

# Update selection to D1
$ws.Range("D1").Select()
